# MAPA_conversions.xlsx — re-ran the unit-merge script, which appended one
# more mapped-unit row ("C" -> "Cz") to the lookup table on Sheet1, and left
# the selection on the cell the user was last working in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended right after the existing A2:D4 block.
$ws.Range("A5").Value = "C"
$ws.Range("B5").Value = "Cz"

# Reflect where the user's cursor ended up after the edit.
$ws.Range("C11").Select()
